$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the file size (D6) and interval duration in minutes (E6).
# Dependent formula cells (F6, J6, K6, L6, M6) recalc automatically.
$ws.Range("D6").Value = 150000
$ws.Range("E6").Value = 5

# Move the active selection from D6 to E6, matching the saved view state.
$ws.Range("E6").Select()
